# Add a new sales record (row 7) to the Planilha1 worksheet, matching the
# data-entry order the author used (vendor name, then date, then CPF,
# then the remaining columns left-to-right), followed by a blank
# formatted row (row 8) carrying the new currency number format forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 -------------------------------------------------------------

# B7: Vendedor (salesperson name) - typed first.
$ws.Range("B7").Value = "Gustavo Sena"

# A7: Data_da_venda - stored as text (matches the existing date cells in
# column A, which use a centered text number format) rather than letting
# Excel auto-convert it to a date serial.
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "11/04/2024"

# C7: CPF_Vendedor
$ws.Range("C7").Value = "333.333.333-34"

# D7: Produto
$ws.Range("D7").Value = "Produto A"

# E7: ID_Produto
$ws.Range("E7").Value = 2

# F7: Cliente
$ws.Range("F7").Value = "Pessoa2"

# G7: CNPJ_CPF_Cliente
$ws.Range("G7").Value = "22.222.222/2222-22"

# H7: Segmento_do_Cliente
$ws.Range("H7").Value = "Segmento B"

# I7: Valor_de_Venda - new currency format "R$ #,##0.00" (no red-negative
# variant, unlike the existing numFmtId 6/8 styles already in the sheet).
$ws.Range("I7").NumberFormat = """R$""\ #,##0.00"
$ws.Range("I7").Value = 1000

# J7: Forma_de_Pagamento
$ws.Range("J7").Value = "Parcelado"

# --- Row 8 ---------------------------------------------------------------
# Blank row below, carrying the new currency format into I8 (no value).
$ws.Range("I8").NumberFormat = """R$""\ #,##0.00"

# --- Selection -------------------------------------------------------
# Leave the selection where the author ended up after the edit.
$ws.Range("B21").Select()
